$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Columns.Item(2).Insert()
$ws.Range("B1").Value = "Execution time (ms)"
$ws.Range("I1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1:B2").Merge()
